$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Data Type"
$ws.Range("A1").Value = "Variable"
$ws.Range("C1").Value = "Values"
$ws.Range("D1").Value = "Notes"

$ws.Range("D7").Select()
